# Test import certificate of Admin
# Update the student roll-number / name / nationality columns and a couple of
# content cells on Sheet1, then leave the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (student #1)
$ws.Range("B2").Value = "HE 130576"
$ws.Range("C2").Value = "Phạm Thanh Hà"
$ws.Range("G2").Value = "Việt Nam"

# Row 3 (student #2)
$ws.Range("B3").Value = "HE 130576@%"
$ws.Range("C3").Value = "Phạm Thanh Hà"
$ws.Range("H3").Value = "Hà Nội"
$ws.Range("K3").ClearContents()

# Row 4 (student #3)
$ws.Range("B4").Value = "HE130576"
$ws.Range("C4").Value = "Phạm Thanh Hà"
$ws.Range("H4").Value = "Hà Nội"
$ws.Range("K4").Value = "Sáo 21"

# Row 5 (student #4)
$ws.Range("C5").Value = "Phạm Thanh Hà"
$ws.Range("G5").Value = "Việt Nam"
$ws.Range("K5").Value = "Lừa trái tim đàn bà @$"

# Match the author's final selection in the saved file
$ws.Range("K11").Select()
